# Update the "Events" sheet: add an "Event" column, mark every row's
# status as "Completed" and track project sub_stage_id / stage_id.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# 1. Insert a new (blank) column before column C. This shifts the
#    existing "Status" column (C) to D and "Date" column (D) to E. The
#    new blank column C inherits the style of the column it pushed.
$ws.Columns("C").Insert()

# 2. Add the new "sub_stage_id" column in F, filled top to bottom.
$ws.Range("F1").Value = "sub_stage_id"
$ws.Range("F2").Value = "A1"
$ws.Range("F3").Value = "B1"
$ws.Range("F4").Value = "C1"
$ws.Range("F5").Value = "A1"
$ws.Range("F6").Value = "B1"

# 3. New column C becomes "Event" and takes the values that used to be
#    in the Status column (now shifted to D).
$ws.Range("C1").Value = "Event"
$ws.Range("C2").Value = $ws.Range("D2").Value2
$ws.Range("C3").Value = $ws.Range("D3").Value2
$ws.Range("C4").Value = $ws.Range("D4").Value2
$ws.Range("C5").Value = $ws.Range("D5").Value2
$ws.Range("C6").Value = $ws.Range("D6").Value2

# 4. The old Status column (now D) is "Completed" for every row.
$ws.Range("D2").Value = "Completed"
$ws.Range("D3").Value = "Completed"
$ws.Range("D4").Value = "Completed"
$ws.Range("D5").Value = "Completed"
$ws.Range("D6").Value = "Completed"

# 5. Add the new "stage_id" column in G, filled top to bottom.
$ws.Range("G1").Value = "stage_id"
$ws.Range("G2").Value = "A"
$ws.Range("G3").Value = "B"
$ws.Range("G4").Value = "C"
$ws.Range("G5").Value = "A"
$ws.Range("G6").Value = "B"

# 6. Copy formatting (style) from existing formatted cells onto the new
#    F:G columns so they match the rest of the table instead of using
#    the workbook's default style.
$ws.Range("A1").Copy()
$ws.Range("F1:G6").PasteSpecial(-4122)

# 7. Row 7 used to hold the "Completed" event for Project DEF; that
#    event no longer applies, so its data is cleared (formatting kept).
#    First extend formatting from row 6 into F7 (which never received a
#    value), then clear the row's contents.
$ws.Range("F6").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("A7:F7").ClearContents()

$excel.CutCopyMode = 0

# 8. Set the column widths to match and restore the active selection.
$ws.Columns("C:C").ColumnWidth = 16.166666666666668
$ws.Columns("D:D").ColumnWidth = 16.166666666666668
$ws.Columns("F:F").ColumnWidth = 11.166666666666666

$ws.Range("G6").Select()
